# Insert a new weekly record at row 147 for "Pepino dulce" (Vega Modelo de
# Temuco), pushing the existing rows 147-191 down to 148-192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 147..191 down to 148..192, leaving a blank row 147 that
# inherits the formatting (incl. the date number format on column D) of
# the row it is inserted above.
$ws.Rows(147).Insert()

# Fill in the newly inserted row with the new data record.
$ws.Cells.Item(147, 1).Value  = 10
$ws.Cells.Item(147, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(147, 3).Value  = "La Araucanía"
$ws.Cells.Item(147, 4).Value  = 44627
$ws.Cells.Item(147, 5).Value  = 9
$ws.Cells.Item(147, 6).Value  = 100112043
$ws.Cells.Item(147, 7).Value  = "Pepino dulce"
$ws.Cells.Item(147, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(147, 9).Value  = "Primera"
$ws.Cells.Item(147, 10).Value = 80
$ws.Cells.Item(147, 11).Value = 12000
$ws.Cells.Item(147, 12).Value = 14000
$ws.Cells.Item(147, 13).Value = 12875
$ws.Cells.Item(147, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(147, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(147, 16).Value = 715
$ws.Cells.Item(147, 17).Value = 18
$ws.Cells.Item(147, 18).Value = "Hortaliza"
